$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1458.5834
$ws.Range("I39").Value = 750.8
$ws.Range("J39").Value = 4997.5
$ws.Range("K39").Value = 2252.4
$ws.Range("L39").Value = 14992.5
$ws.Range("M39").Value = -1956.4
$ws.Range("N39").Value = -15584.5

$ws.Range("H137").Value = 19121902
$ws.Range("I137").Value = 38463436
$ws.Range("J137").Value = 1161904.6
$ws.Range("K137").Value = 115390308
$ws.Range("L137").Value = 3485713.8
$ws.Range("M137").Value = -115387758
$ws.Range("N137").Value = -3490813.8

$ws.Range("H138").Value = 3064.6584
$ws.Range("I138").Value = 1386.0435
$ws.Range("J138").Value = 5209.5557
$ws.Range("K138").Value = 4158.1305
$ws.Range("L138").Value = 15628.6671
$ws.Range("M138").Value = 981.8694999999998
$ws.Range("N138").Value = -25908.6671

$ws.Range("H141").Value = 1958.0667
$ws.Range("I141").Value = 1976.6428
$ws.Range("J141").Value = 1698
$ws.Range("K141").Value = 5929.928400000001
$ws.Range("L141").Value = 5094
$ws.Range("M141").Value = -749.9284000000007
$ws.Range("N141").Value = -15454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38040.867
$ws.Range("I32").Value = 40735.18
$ws.Range("J32").Value = 320.5
$ws.Range("K32").Value = 40735.18
$ws.Range("L32").Value = 320.5
$ws.Range("M32").Value = -40448.18
$ws.Range("N32").Value = -894.5

$ws.Range("H45").Value = 2468.16
$ws.Range("I45").Value = 1747.8096
$ws.Range("J45").Value = 6250
$ws.Range("K45").Value = 1747.8096
$ws.Range("L45").Value = 6250
$ws.Range("M45").Value = -1370.8096
$ws.Range("N45").Value = -7004

$ws.Range("H55").Value = 99978
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 99978
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 99978
$ws.Range("N55").Value = -100608

$ws.Range("H74").Value = 3313.9285
$ws.Range("I74").Value = 932.7778
$ws.Range("J74").Value = 7600
$ws.Range("K74").Value = 932.7778
$ws.Range("L74").Value = 7600
$ws.Range("M74").Value = -58.77779999999996
$ws.Range("N74").Value = -9348

$ws.Range("H77").Value = 3313.9285
$ws.Range("I77").Value = 932.7778
$ws.Range("J77").Value = 7600
$ws.Range("K77").Value = 4663.889
$ws.Range("L77").Value = 38000
$ws.Range("M77").Value = -295.8890000000001
$ws.Range("N77").Value = -46736

$ws.Range("H113").Value = 88948.875
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 88948.875
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 88948.875
$ws.Range("N113").Value = -97626.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1389.4546
$ws.Range("I20").Value = 1310.875
$ws.Range("J20").Value = 1599
$ws.Range("K20").Value = 1310.875
$ws.Range("L20").Value = 1599
$ws.Range("M20").Value = -1063.875
$ws.Range("N20").Value = -2093

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3575806.8
$ws.Range("I4").Value = 7151571
$ws.Range("J4").Value = 42.857143
$ws.Range("K4").Value = 7151571
$ws.Range("L4").Value = 42.857143
$ws.Range("M4").Value = -7151459
$ws.Range("N4").Value = -266.857143

$ws.Range("H58").Value = 592056.4399999999
$ws.Range("I58").Value = 825884
$ws.Range("J58").Value = 7487.5
$ws.Range("K58").Value = 825884
$ws.Range("L58").Value = 7487.5
$ws.Range("M58").Value = -825681
$ws.Range("N58").Value = -7893.5

$ws.Range("H95").Value = 17666.334
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 17666.334
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 17666.334
$ws.Range("N95").Value = -23158.334

$ws.Range("H96").Value = 5760.697
$ws.Range("I96").Value = 4500
$ws.Range("J96").Value = 5934.5864
$ws.Range("K96").Value = 4500
$ws.Range("L96").Value = 5934.5864
$ws.Range("M96").Value = -1754
$ws.Range("N96").Value = -11426.5864

$ws.Range("H99").Value = 3499.6667
$ws.Range("I99").Value = 3074.5
$ws.Range("J99").Value = 4350
$ws.Range("K99").Value = 3074.5
$ws.Range("L99").Value = 4350
$ws.Range("M99").Value = -1576.5
$ws.Range("N99").Value = -7346

$ws.Range("H126").Value = 3499.6667
$ws.Range("I126").Value = 3074.5
$ws.Range("J126").Value = 4350
$ws.Range("K126").Value = 9223.5
$ws.Range("L126").Value = 13050
$ws.Range("M126").Value = -6753.5
$ws.Range("N126").Value = -17990

$ws.Range("H134").Value = 7866.773
$ws.Range("I134").Value = 8074.7144
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 24224.1432
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -21689.1432
$ws.Range("N134").Value = -15570

$ws.Range("H136").Value = 592056.4399999999
$ws.Range("I136").Value = 825884
$ws.Range("J136").Value = 7487.5
$ws.Range("K136").Value = 2477652
$ws.Range("L136").Value = 22462.5
$ws.Range("M136").Value = -2475102
$ws.Range("N136").Value = -27562.5

$ws.Range("H140").Value = 95000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 95000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 95000
$ws.Range("N140").Value = -105360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 2844.5
$ws.Range("I45").Value = 3000
$ws.Range("J45").Value = 2792.6667
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 8378.000100000001
$ws.Range("M45").Value = -8468
$ws.Range("N45").Value = -9442.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H41").Value = 4739.25
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 5652.3335
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 5652.3335
$ws.Range("M41").Value = -1645
$ws.Range("N41").Value = -6362.3335

$ws.Range("H42").Value = 80000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 80000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 80000
$ws.Range("N42").Value = -80970

$ws.Range("H102").Value = 3091.282
$ws.Range("I102").Value = 2157.2173
$ws.Range("J102").Value = 4434
$ws.Range("K102").Value = 2157.2173
$ws.Range("L102").Value = 4434
$ws.Range("M102").Value = -535.2172999999998
$ws.Range("N102").Value = -7678

$ws.Range("H115").Value = 80000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 80000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 80000
$ws.Range("N115").Value = -82350

$ws.Range("H122").Value = 4976.5713
$ws.Range("I122").Value = 2758.0625
$ws.Range("J122").Value = 7934.5835
$ws.Range("K122").Value = 8274.1875
$ws.Range("L122").Value = 23803.7505
$ws.Range("M122").Value = -5824.1875
$ws.Range("N122").Value = -28703.7505

$ws.Range("H132").Value = 56230784
$ws.Range("I132").Value = 72293200
$ws.Range("J132").Value = 12328.25
$ws.Range("K132").Value = 216879600
$ws.Range("L132").Value = 36984.75
$ws.Range("M132").Value = -216877070
$ws.Range("N132").Value = -42044.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1303.2222
$ws.Range("I22").Value = 807.25
$ws.Range("J22").Value = 1700
$ws.Range("K22").Value = 807.25
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = -512.25
$ws.Range("N22").Value = -2290

$ws.Range("H27").Value = 1303.2222
$ws.Range("I27").Value = 807.25
$ws.Range("J27").Value = 1700
$ws.Range("K27").Value = 807.25
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = -700.25
$ws.Range("N27").Value = -1914

$ws.Range("H55").Value = 180.6923
$ws.Range("I55").Value = 160.66667
$ws.Range("J55").Value = 197.85715
$ws.Range("K55").Value = 160.66667
$ws.Range("L55").Value = 197.85715
$ws.Range("M55").Value = 12.33332999999999
$ws.Range("N55").Value = -543.85715

$ws.Range("H68").Value = 2999
$ws.Range("I68").Value = 2999
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2999
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2250

$ws.Range("H71").Value = 2999
$ws.Range("I71").Value = 2999
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14995
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -11251

$ws.Range("H82").Value = 946.75
$ws.Range("I82").Value = 844.8333
$ws.Range("J82").Value = 1252.5
$ws.Range("K82").Value = 844.8333
$ws.Range("L82").Value = 1252.5
$ws.Range("M82").Value = -483.8333
$ws.Range("N82").Value = -1974.5

$ws.Range("H85").Value = 946.75
$ws.Range("I85").Value = 844.8333
$ws.Range("J85").Value = 1252.5
$ws.Range("K85").Value = 844.8333
$ws.Range("L85").Value = 1252.5
$ws.Range("M85").Value = 403.1667
$ws.Range("N85").Value = -3748.5

$ws.Range("H122").Value = 4357.905
$ws.Range("I122").Value = 3789.8708
$ws.Range("J122").Value = 5958.727
$ws.Range("K122").Value = 11369.6124
$ws.Range("L122").Value = 17876.181
$ws.Range("M122").Value = -8919.6124
$ws.Range("N122").Value = -22776.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3212.125
$ws.Range("I81").Value = 2033
$ws.Range("J81").Value = 6749.5
$ws.Range("K81").Value = 4066
$ws.Range("L81").Value = 13499
$ws.Range("M81").Value = -3005
$ws.Range("N81").Value = -15621

$ws.Range("H84").Value = 3212.125
$ws.Range("I84").Value = 2033
$ws.Range("J84").Value = 6749.5
$ws.Range("K84").Value = 20330
$ws.Range("L84").Value = 67495
$ws.Range("M84").Value = -15026
$ws.Range("N84").Value = -78103
